$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Value = 1.35
$ws.Range("R2").Value = 1.46
$ws.Range("J3").Value = 3.25
$ws.Range("X3").Value = 21
$ws.Range("AJ3").Value = 32
$ws.Range("Y4").Value = 980
$ws.Range("M5").Value = 1.09
$ws.Range("N5").Value = 2.94
$ws.Range("R5").Value = 1.24
$ws.Range("S5").Value = 4
$ws.Range("T5").Value = 1.85
$ws.Range("U5").Value = 1.89
$ws.Range("X5").Value = 12
$ws.Range("Y5").Value = 13
$ws.Range("Z5").Value = 28
$ws.Range("AA5").Value = 85
$ws.Range("AB5").Value = 9.2
$ws.Range("AC5").Value = 8
$ws.Range("AD5").Value = 17
$ws.Range("AE5").Value = 55
$ws.Range("AF5").Value = 15
$ws.Range("AG5").Value = 12.5
$ws.Range("AH5").Value = 21
$ws.Range("AI5").Value = 70
$ws.Range("AJ5").Value = 980
$ws.Range("AK5").Value = 980
$ws.Range("AL5").Value = 50
$ws.Range("AN5").Value = 27
$ws.Range("F6").Value = 1.83
$ws.Range("H6").Value = 4.4
$ws.Range("I6").Value = 5.7
$ws.Range("L6").Value = 1.42
$ws.Range("M6").Value = 1.1
$ws.Range("N6").Value = 2.78
$ws.Range("O6").Value = 1.45
$ws.Range("Q6").Value = 2.22
$ws.Range("R6").Value = 1.22
$ws.Range("S6").Value = 4.5
$ws.Range("T6").Value = 2.06
$ws.Range("U6").Value = 1.74
$ws.Range("V6").Value = 1.21
$ws.Range("X6").Value = 1000
$ws.Range("AB6").Value = 8.4
$ws.Range("AC6").Value = 9.8
$ws.Range("AJ6").Value = 26
$ws.Range("AK6").Value = 30
$ws.Range("AN6").Value = 22
$ws.Range("J7").Value = 3.6
$ws.Range("Y7").Value = 1000
$ws.Range("AJ7").Value = 1000
$ws.Range("L8").Value = 1.46
$ws.Range("M8").Value = 1.09
$ws.Range("N8").Value = 3.1
$ws.Range("Q8").Value = 2.06
$ws.Range("R8").Value = 1.27
$ws.Range("S8").Value = 3.95
$ws.Range("T8").Value = 1.84
$ws.Range("U8").Value = 1.98
$ws.Range("W8").Value = 1.65
$ws.Range("X8").Value = 12
$ws.Range("Y8").Value = 12
$ws.Range("Z8").Value = 24
$ws.Range("AA8").Value = 70
$ws.Range("AB8").Value = 9.4
$ws.Range("AC8").Value = 7.8
$ws.Range("AD8").Value = 15.5
$ws.Range("AE8").Value = 980
$ws.Range("AF8").Value = 15.5
$ws.Range("AG8").Value = 12
$ws.Range("AH8").Value = 19.5
$ws.Range("AI8").Value = 60
$ws.Range("AJ8").Value = 36
$ws.Range("AK8").Value = 30
$ws.Range("AL8").Value = 48
$ws.Range("AM8").Value = 140
$ws.Range("AN8").Value = 32
$ws.Range("AO8").Value = 60
$ws.Range("N9").Value = 1.74
$ws.Range("P9").Value = 1.73
$ws.Range("Q9").Value = 1.95
